# Apply the seat-chart update for ICC S2C5 (21 Oct, 2022)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Seats")

# Header / date
$ws.Range("A1").Value = "ICC S2C5 座位表"
$ws.Range("K1").Value = " 21 Oct, 2022"

# Row 4
$ws.Range("B4").Value = "廖从云"
$ws.Range("C4").Value = "边麓元"
$ws.Range("E4").Value = "陈元畅"
$ws.Range("F4").Value = "詹悦"
$ws.Range("H4").Value = "李星宸"
$ws.Range("I4").Value = "龙飞宇"
$ws.Range("K4").Value = "石清泓"

# Row 5
$ws.Range("B5").Value = "王昊天"
$ws.Range("C5").Value = "丁鹏元"
$ws.Range("E5").Value = "张宸瑞"
$ws.Range("F5").Value = "曾韦翔"
$ws.Range("H5").Value = "龚搏扬"
$ws.Range("I5").Value = "程启航"
$ws.Range("K5").Value = "郑俊永"
$ws.Range("L5").Value = ""

# Row 6
$ws.Range("B6").Value = "林彦含"
$ws.Range("E6").Value = "范青桐"
$ws.Range("F6").Value = "邱晨朔"
$ws.Range("H6").Value = "杨熙宇"
$ws.Range("I6").Value = "张扬"
$ws.Range("K6").Value = "迟涵予"
$ws.Range("L6").Value = "白宇轩"

# Row 7
$ws.Range("B7").Value = "吴周毅"
$ws.Range("C7").Value = "卢逸"
$ws.Range("K7").Value = "蔡朋骏"
$ws.Range("L7").Value = "陈李石农"

# Update the active selection to A2 (matches the saved view state)
$ws.Range("A2").Select()
